$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# 1. The shared status text "Ready for handoff" becomes "Handback transform failed"
#    for the 7f4aa704 file row, everywhere it is used (Overview + per-locale sheets).
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# 2. Populate the "Error Detail" column (P) for the 7f4aa704 file row on each
#    locale sheet with the handback/handoff file name mismatch message.
$wsZhCn.Range("P3").Value = "Handback file name: puglxqu5.bnc is different with handoff file name: 7f4aa704-d9eb-4e4b-8a1f-fe2f21e56be6.12c4a5e617681cac84ad8e9692c985f2b3e99990.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: puglxqu5.bnc is different with handoff file name: 7f4aa704-d9eb-4e4b-8a1f-fe2f21e56be6.12c4a5e617681cac84ad8e9692c985f2b3e99990.de-de."

# 3. Widen the "Error Detail" column (16 / P) to fit the new messages.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
